# Updated cryptos list values (Price / Volume(1h)) for the rows that changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number need an explicit
# Text format first, otherwise Excel auto-converts the literal into a floating
# point number (losing the original "40.15" style text representation).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "40.143.06"
$ws.Range("E2").Value = "  -2.15%  "

$ws.Range("D3").Value = "2.348.34"
$ws.Range("E3").Value = "  -3.03%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "311.03"
$ws.Range("E5").Value = "  -1.60%  "

$ws.Range("D6").Value = "85.66"
$ws.Range("E6").Value = "  -3.46%  "

$ws.Range("E7").Value = "  -1.50%  "

$ws.Range("D9").Value = "0.484"
$ws.Range("E9").Value = "  -2.10%  "

$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  -2.18%  "

$ws.Range("D11").Value = "30.13"
$ws.Range("E11").Value = "  -5.51%  "

$ws.Range("E12").Value = "  +1.28%  "

$ws.Range("D13").Value = "2.710.85"
$ws.Range("E13").Value = "  -3.21%  "

$ws.Range("E14").Value = "  -3.63%  "

$ws.Range("D15").Value = "14.81"
$ws.Range("E15").Value = "  -4.79%  "

$ws.Range("D16").Value = "2.372.26"
$ws.Range("E16").Value = "  -1.94%  "

$ws.Range("E17").Value = "  -0.88%  "

$ws.Range("D18").Value = "40.127.22"
$ws.Range("E18").Value = "  -2.06%  "

$ws.Range("D19").Value = "0.0₃0903"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").Value = "6.11"
$ws.Range("E20").Value = "  -1.90%  "

$ws.Range("D21").Value = "68.18"
$ws.Range("E21").Value = "  -5.01%  "

$ws.Range("D22").Value = "10.73"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("D23").Value = "235.35"
$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("E24").Value = "  -5.00%  "

$ws.Range("E26").Value = "  -2.03%  "

$ws.Range("D27").Value = "23.61"
$ws.Range("E27").Value = "  -1.75%  "

$ws.Range("E28").Value = "  -3.19%  "

$ws.Range("D29").Value = "9.27"
$ws.Range("E29").Value = "  -2.54%  "

$ws.Range("D30").Value = "34.70"
$ws.Range("E30").Value = "  +0.08%  "

$ws.Range("D31").Value = "153.89"
$ws.Range("E31").Value = "  -1.24%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("D33").Value = "5.12"
$ws.Range("E33").Value = "  -2.23%  "

$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("D35").Value = "0.0722"
$ws.Range("E35").Value = "  -2.70%  "

$ws.Range("E36").Value = "  -0.57%  "

$ws.Range("E37").Value = "  -3.79%  "

$ws.Range("D38").Value = "0.0984"
$ws.Range("E38").Value = "  -1.52%  "

$ws.Range("D39").Value = "15.62"
$ws.Range("E39").Value = "  -6.12%  "

$ws.Range("D40").Value = "1.72"
$ws.Range("E40").Value = "  -2.29%  "

$ws.Range("E41").Value = "  +0.74%  "

$ws.Range("D42").Value = "1.963.56"
$ws.Range("E42").Value = "  -0.99%  "

$ws.Range("E43").Value = "  -3.02%  "

$ws.Range("E44").Value = "  -3.40%  "

$ws.Range("D45").Value = "17.68"
$ws.Range("E45").Value = "  -4.61%  "

$ws.Range("D46").Value = "9.37"
$ws.Range("E46").Value = "  -0.68%  "

$ws.Range("D47").Value = "2.71"
$ws.Range("E47").Value = "  -5.58%  "

$ws.Range("D48").Value = "2.569.81"
$ws.Range("E48").Value = "  -3.44%  "

$ws.Range("D49").Value = "93.32"
$ws.Range("E49").Value = "  -1.20%  "

$ws.Range("D50").Value = "70.68"
$ws.Range("E50").Value = "  -2.99%  "

$ws.Range("D51").Value = "50.55"
$ws.Range("E51").Value = "  -2.18%  "
